$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the hard-coded Total (column B) values for rows 2-31 with a formula
# that sums the sector columns (C:I) plus the two extra columns (K, L).
# B2 is set on its own (ordinary formula), and B3:B31 are set together so
# Excel groups them as one shared formula, matching how the workbook was
# originally authored (fill down from B3).
$ws.Range("B2").Formula = "=SUM(C2:I2)+K2+L2"
$ws.Range("B3:B31").Formula = "=SUM(C3:I3)+K3+L3"

# Update the active selection to match the saved state in the workbook.
$ws.Range("D19").Select()
